# linux managed systemler ok
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Linux target systems: fix the standby3 host entry and add the bastion
# host entry for the "root" / SSH-managed rows.
$ws.Range("B6").Value = "standby3.quasys.locals"
$ws.Range("B10").Value = "bastion.quasys.local"

# Move the active selection on Sheet1 to B13.
$ws.Activate()
$ws.Range("B13").Select()

# Reflect the user's final window arrangement.
$w = $excel.ActiveWindow
$w.Left = 38280
$w.Top = -240
$w.Width = 29040
$w.Height = 15720
